# Adding pickling to replication
# Fill in the standard-error rows (theta_se = row 4, lambda_se = row 6)
# which were previously placeholder "(nan)" text values, with the actual
# computed standard errors (as text, matching the "(x.xx)" display style
# used throughout this table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$thetaSe = [ordered]@{
    "B4" = "(8.64)"
    "C4" = "(1.31)"
    "D4" = "(0.5)"
    "E4" = "(1.92)"
    "F4" = "(0.96)"
    "G4" = "(1.41)"
    "H4" = "(3.23)"
    "I4" = "(2.16)"
    "J4" = "(0.22)"
    "K4" = "(4.27)"
    "L4" = "(3.48)"
}

$lambdaSe = [ordered]@{
    "B6" = "(4.87)"
    "C6" = "(1.08)"
    "D6" = "(0.37)"
    "E6" = "(1.63)"
    "F6" = "(0.28)"
    "G6" = "(0.56)"
    "H6" = "(2.09)"
    "I6" = "(3.99)"
    "J6" = "(0.29)"
    "K6" = "(3.64)"
    "L6" = "(3.43)"
}

foreach ($addr in $thetaSe.Keys) {
    $ws.Range($addr).Value = $thetaSe[$addr]
}

foreach ($addr in $lambdaSe.Keys) {
    $ws.Range($addr).Value = $lambdaSe[$addr]
}
